# Update the answer cells in the three-digit ÷ one-digit division table
# to the newly generated problem set (quotient, remainder pairs).
$d = $word.ActiveDocument
$d.Content.Find.Execute("250÷5=50, 0", $true, $false, $false, $false, $false, $true, 1, $false, "105÷9=11, 6", 2) | Out-Null
$d.Content.Find.Execute("977÷4=244, 1", $true, $false, $false, $false, $false, $true, 1, $false, "847÷9=94, 1", 2) | Out-Null
$d.Content.Find.Execute("931÷8=116, 3", $true, $false, $false, $false, $false, $true, 1, $false, "900÷7=128, 4", 2) | Out-Null
$d.Content.Find.Execute("446÷7=63, 5", $true, $false, $false, $false, $false, $true, 1, $false, "954÷8=119, 2", 2) | Out-Null
$d.Content.Find.Execute("380÷5=76, 0", $true, $false, $false, $false, $false, $true, 1, $false, "514÷9=57, 1", 2) | Out-Null
$d.Content.Find.Execute("818÷9=90, 8", $true, $false, $false, $false, $false, $true, 1, $false, "733÷7=104, 5", 2) | Out-Null
$d.Content.Find.Execute("715÷5=143, 0", $true, $false, $false, $false, $false, $true, 1, $false, "871÷4=217, 3", 2) | Out-Null
$d.Content.Find.Execute("332÷7=47, 3", $true, $false, $false, $false, $false, $true, 1, $false, "427÷3=142, 1", 2) | Out-Null
$d.Content.Find.Execute("496÷6=82, 4", $true, $false, $false, $false, $false, $true, 1, $false, "504÷6=84, 0", 2) | Out-Null
$d.Content.Find.Execute("492÷7=70, 2", $true, $false, $false, $false, $false, $true, 1, $false, "427÷2=213, 1", 2) | Out-Null
$d.Content.Find.Execute("293÷4=73, 1", $true, $false, $false, $false, $false, $true, 1, $false, "627÷2=313, 1", 2) | Out-Null
$d.Content.Find.Execute("960÷2=480, 0", $true, $false, $false, $false, $false, $true, 1, $false, "259÷7=37, 0", 2) | Out-Null
$d.Content.Find.Execute("749÷5=149, 4", $true, $false, $false, $false, $false, $true, 1, $false, "706÷4=176, 2", 2) | Out-Null
$d.Content.Find.Execute("192÷9=21, 3", $true, $false, $false, $false, $false, $true, 1, $false, "471÷3=157, 0", 2) | Out-Null
$d.Content.Find.Execute("133÷5=26, 3", $true, $false, $false, $false, $false, $true, 1, $false, "209÷8=26, 1", 2) | Out-Null
$d.Content.Find.Execute("181÷7=25, 6", $true, $false, $false, $false, $false, $true, 1, $false, "222÷6=37, 0", 2) | Out-Null
$d.Content.Find.Execute("571÷3=190, 1", $true, $false, $false, $false, $false, $true, 1, $false, "905÷3=301, 2", 2) | Out-Null
$d.Content.Find.Execute("165÷3=55, 0", $true, $false, $false, $false, $false, $true, 1, $false, "505÷8=63, 1", 2) | Out-Null
$d.Content.Find.Execute("144÷4=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "909÷7=129, 6", 2) | Out-Null
$d.Content.Find.Execute("463÷9=51, 4", $true, $false, $false, $false, $false, $true, 1, $false, "476÷6=79, 2", 2) | Out-Null
$d.Content.Find.Execute("325÷5=65, 0", $true, $false, $false, $false, $false, $true, 1, $false, "620÷6=103, 2", 2) | Out-Null
$d.Content.Find.Execute("271÷7=38, 5", $true, $false, $false, $false, $false, $true, 1, $false, "253÷9=28, 1", 2) | Out-Null
$d.Content.Find.Execute("866÷5=173, 1", $true, $false, $false, $false, $false, $true, 1, $false, "367÷7=52, 3", 2) | Out-Null
$d.Content.Find.Execute("699÷8=87, 3", $true, $false, $false, $false, $false, $true, 1, $false, "829÷2=414, 1", 2) | Out-Null
$d.Content.Find.Execute("853÷6=142, 1", $true, $false, $false, $false, $false, $true, 1, $false, "873÷4=218, 1", 2) | Out-Null
